$p = $ppt.ActivePresentation
try {
  $p.SaveAs("C:\evals\00000000-0000-0000-0000-000000000000\before.pptx")
  Write-Output "SaveAs worked"
} catch {
  Write-Output "SaveAs err: $_"
}
